$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A64").Value = "Journal of Politics in Latin America"
$ws.Range("B64").Value = "<a href='https://journals.sagepub.com/author-instructions/PLA'target='_blank'>Research Note</a>"
$ws.Range("C64").Value = "3k -- 5k words"
$ws.Range("D64").Value = "*N/A*"

$ws.Range("C64").Select()
